$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.529.16"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "2.470.54"
$ws.Range("E3").Value = "  -0.68%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.98"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "92.03"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.549"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  +3.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.53"
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").Value = "2.849.92"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.86"
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.99"
$ws.Range("E15").Value = "  +3.40%  "
$ws.Range("D16").Value = "2.474.19"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.778"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "41.570.19"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.16"
$ws.Range("E21").Value = "  +2.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.13"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.22"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("E25").Value = "  +0.65%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.83"
$ws.Range("E27").Value = "  +2.91%  "
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.70"
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.46"
$ws.Range("E30").Value = "  -4.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.16"
$ws.Range("E31").Value = "  +1.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.46"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.58"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.24"
$ws.Range("E35").Value = "  -3.85%  "
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.89"
$ws.Range("E37").Value = "  -6.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.104"
$ws.Range("E38").Value = "  +2.46%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.79"
$ws.Range("E40").Value = "  -4.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("E41").Value = "  -2.56%  "
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").Value = "1.946.73"
$ws.Range("E43").Value = "  -1.98%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.98"
$ws.Range("E44").Value = "  -3.70%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0283"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.09"
$ws.Range("E47").Value = "  +2.45%  "
$ws.Range("D48").Value = "2.707.95"
$ws.Range("E48").Value = "  -0.86%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.44"
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "67.26"
$ws.Range("E50").Value = "  -3.11%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.70"
$ws.Range("E51").Value = "  +3.65%  "
